# Apply the 11.b.2 worksheet update:
#  - extend the yearly data table from 2019-only to 2019-2023
#  - D4 becomes a genuine number (484) instead of a text value "484"
#  - header row gets taller, data columns A:C get slightly narrower
#  - selection resets to the default cell (A1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 header: make it taller to fit the (now wider) header text ---
$ws.Rows(1).RowHeight = 79.5

# --- Column widths for A:C (41.28515625 -> 40.140625 characters) ---
$ws.Columns("A:C").ColumnWidth = 39.33

# --- Fix D4: was stored as text "484", must become numeric 484 ---
$ws.Range("D4").Value = 484

# --- Row 3: year headers for the new columns, matching D3's format ---
$ws.Range("D3").Copy()
$ws.Range("E3:H3").PasteSpecial(-4122)
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# --- Row 4: number of local governments, matching D4's format ---
$ws.Range("D4").Copy()
$ws.Range("E4:H4").PasteSpecial(-4122)
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# --- Row 5: proportion (%), matching D5's format ---
$ws.Range("D5").Copy()
$ws.Range("E5:H5").PasteSpecial(-4122)
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# --- Row 6: count of local governments adopting DRR strategies, matching D6's format ---
$ws.Range("D6").Copy()
$ws.Range("E6:H6").PasteSpecial(-4122)
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# --- Reset the active selection back to the default A1 cell ---
$ws.Range("A1").Select() | Out-Null

$excel.CutCopyMode = $false
